# SAV-700: Update charts importer test fixture
# - Ensure ID (code) of all complex chart core questions is constant:
#   the "code" column (A) on the "Core" sheet is set equal to the
#   question's "type" column (B) value instead of the old
#   testchartcorecode0..3 placeholders.
# - Ensure ID of the charting date recorded program data element is
#   constant: "testchartcode0" -> "PatientChartingDate" on the
#   "Test Chart" sheet.
# - The cell-highlight style that used to sit on the Metadata sheet's
#   isSensitive value (G8) now sits on the "Test Chart" sheet's
#   PatientChartingDate code cell (A2) instead.

$wb = $excel.ActiveWorkbook

$core = $wb.Worksheets.Item("Core")
$testChart = $wb.Worksheets.Item("Test Chart")
$metadata = $wb.Worksheets.Item("Metadata")

# --- Move the distinctive cell style from Metadata!G8 to "Test Chart"!A2
# before the values change, so the style follows the right cell.
$metadata.Range("G8").Copy()
$testChart.Range("A2").PasteSpecial(-4122)

$metadata.Range("G7").Copy()
$metadata.Range("G8").PasteSpecial(-4122)

# --- Core sheet: make the "code" column equal to the "type" column
# for each of the complex chart core questions.
$core.Range("A2").Value = $core.Range("B2").Value2
$core.Range("A3").Value = $core.Range("B3").Value2
$core.Range("A4").Value = $core.Range("B4").Value2
$core.Range("A5").Value = $core.Range("B5").Value2

# --- Test Chart sheet: rename the charting-date question code.
$testChart.Range("A2").Value = "PatientChartingDate"
